$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Reading list")
$ws.Activate()

# The paper that has now been read through gets removed from the reading
# list entirely - delete its whole row (row 6), shifting the rows below up.
$ws.Rows.Item(6).Delete()

# Mark the rows that have now been read (rows 2-8 in column B) with a
# green fill, matching the "read"/processed colour used elsewhere in the
# workbook.
$readRange = $ws.Range("B2:B8")
$readRange.Interior.Color = 0x50B000

$ws.Range("A16").Select()
